# Insert a new weekly price record as row 5 in the "Bruselas (repollito)"
# Feria Lagunitas de Puerto Montt sheet. Existing rows 5-23 shift down to
# 6-24 (their data is unchanged, only their row position moves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 5..23 down to 6..24, preserving their contents/formatting.
$ws.Rows.Item(5).Insert()

# Populate the newly-opened row 5 with the new record.
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(5, 3).Value = "Los Lagos"
$ws.Cells.Item(5, 4).Value = 44819
$ws.Cells.Item(5, 5).Value = 10
$ws.Cells.Item(5, 6).Value = 100112035
$ws.Cells.Item(5, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 70
$ws.Cells.Item(5, 11).Value = 22000
$ws.Cells.Item(5, 12).Value = 22000
$ws.Cells.Item(5, 13).Value = 22000
$ws.Cells.Item(5, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(5, 16).Value = 1467
$ws.Cells.Item(5, 17).Value = 15
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# Match the date formatting used by the other rows' "Fecha" column.
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat
